$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quarterly indexing bug-fix: a new quarter (Q6 data at row 2) was inserted
# at the top of the rolling window, shifting the previously-computed rows
# down by one. Row 9 (Q13) additionally received a corrected recomputation.

$ws.Cells.Item(2, 2).Value = 0.06673029405116952
$ws.Cells.Item(2, 3).Value = 0.3004351988545212
$ws.Cells.Item(2, 4).Value = 0.1299092234646646
$ws.Cells.Item(2, 5).Value = 0.3604292211581416
$ws.Cells.Item(2, 6).Value = 0.3675687354335792
$ws.Cells.Item(2, 7).Value = 14

$ws.Cells.Item(3, 2).Value = 0.1619789639338207
$ws.Cells.Item(3, 3).Value = 0.3076492955603238
$ws.Cells.Item(3, 4).Value = 0.1317988477936555
$ws.Cells.Item(3, 5).Value = 0.3630411103355314
$ws.Cells.Item(3, 6).Value = 0.3381692893945722
$ws.Cells.Item(3, 7).Value = 13

$ws.Cells.Item(4, 2).Value = 0.1981620169181928
$ws.Cells.Item(4, 3).Value = 0.3116209034719995
$ws.Cells.Item(4, 4).Value = 0.1296824127478638
$ws.Cells.Item(4, 5).Value = 0.3601144439589501
$ws.Cells.Item(4, 6).Value = 0.3140600309705328
$ws.Cells.Item(4, 7).Value = 12

$ws.Cells.Item(5, 2).Value = 0.2201334750512514
$ws.Cells.Item(5, 3).Value = 0.2700445318363022
$ws.Cells.Item(5, 4).Value = 0.09544857498369073
$ws.Cells.Item(5, 5).Value = 0.3089475278808535
$ws.Cells.Item(5, 6).Value = 0.2273517340160526
$ws.Cells.Item(5, 7).Value = 11

$ws.Cells.Item(6, 2).Value = 0.2251408850412157
$ws.Cells.Item(6, 3).Value = 0.2596564006926809
$ws.Cells.Item(6, 4).Value = 0.09977641923863516
$ws.Cells.Item(6, 5).Value = 0.3158740559758512
$ws.Cells.Item(6, 6).Value = 0.2335427658231482
$ws.Cells.Item(6, 7).Value = 10

$ws.Cells.Item(7, 2).Value = 0.1781092598615317
$ws.Cells.Item(7, 3).Value = 0.1781092598615317
$ws.Cells.Item(7, 4).Value = 0.05114199694748629
$ws.Cells.Item(7, 5).Value = 0.2261459638098507
$ws.Cells.Item(7, 6).Value = 0.147805529536099
$ws.Cells.Item(7, 7).Value = 9

$ws.Cells.Item(8, 2).Value = 0.2144878488911046
$ws.Cells.Item(8, 3).Value = 0.2313407205340803
$ws.Cells.Item(8, 4).Value = 0.09456382460910477
$ws.Cells.Item(8, 5).Value = 0.3075123161909207
$ws.Cells.Item(8, 6).Value = 0.235575374136654
$ws.Cells.Item(8, 7).Value = 8

$ws.Cells.Item(9, 2).Value = 0.249316951555495
$ws.Cells.Item(9, 3).Value = 0.2697110612940698
$ws.Cells.Item(9, 4).Value = 0.08937093368650917
$ws.Cells.Item(9, 5).Value = 0.2989497176558445
$ws.Cells.Item(9, 6).Value = 0.1781777855004607
$ws.Cells.Item(9, 7).Value = 7

$ws.Cells.Item(10, 2).Value = 0.2987280035122604
$ws.Cells.Item(10, 3).Value = 0.2987280035122604
$ws.Cells.Item(10, 4).Value = 0.1585806734357395
$ws.Cells.Item(10, 5).Value = 0.3982218896993729
$ws.Cells.Item(10, 6).Value = 0.2844280147574629
$ws.Cells.Item(10, 7).Value = 7

$ws.Cells.Item(11, 2).Value = 0.07168159692063568
$ws.Cells.Item(11, 3).Value = 0.1326419446783583
$ws.Cells.Item(11, 4).Value = 0.02840911087054501
$ws.Cells.Item(11, 5).Value = 0.1685500248310424
$ws.Cells.Item(11, 6).Value = 0.1705537288270639
$ws.Cells.Item(11, 7).Value = 5

